$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "60-38=",
    "49+19=",
    "72+9=",
    "91-58=",
    "61-47=",
    "84-9=",
    "61-27=",
    "93-4=",
    "56+9=",
    "72+19=",
    "61-6=",
    "7+74=",
    "97-39=",
    "35-28=",
    "15+77=",
    "37-9=",
    "8+87=",
    "57+39=",
    "95-89=",
    "84-76=",
    "63-58=",
    "19+34=",
    "7+6=",
    "4+79=",
    "67+18=",
    "77+7=",
    "90-2=",
    "34-18=",
    "57+15=",
    "85-67=",
    "13+49=",
    "55-36=",
    "61-43=",
    "61-28=",
    "58+36=",
    "97-38=",
    "73-19=",
    "8+39=",
    "15+38=",
    "47+9=",
    "62-57=",
    "77-38=",
    "32+19=",
    "68+28=",
    "70-59=",
    "43-4=",
    "66+17=",
    "69+13=",
    "56+29=",
    "98-29=",
    "78+13=",
    "90-52=",
    "19+29=",
    "18+66=",
    "35+59=",
    "43+19=",
    "59+26=",
    "61-13=",
    "17+66=",
    "51-15=",
    "59+6=",
    "4+27=",
    "24+57=",
    "36+57=",
    "69+6=",
    "72-54=",
    "58+8=",
    "36+18=",
    "56+27=",
    "84-76=",
    "28+69=",
    "35-6=",
    "23+18=",
    "38+58=",
    "19+55=",
    "27+27=",
    "15+36=",
    "51-7=",
    "28+63=",
    "61-26=",
    "34+28=",
    "61-37=",
    "90-5=",
    "35+27=",
    "70-52=",
    "18+68=",
    "70-42=",
    "92-27=",
    "29+37=",
    "90-23=",
    "61-46=",
    "48+4=",
    "28+19=",
    "38+37=",
    "95-37=",
    "24+28=",
    "5+88=",
    "36+16=",
    "15+7=",
    "94-18=",
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $expected = $newValues[$idx]
        $cell.Range.Text = $expected
        $idx++
    }
}

Write-Host "Done applying $idx replacements"